# Insert 9 new rows of historical data into the sheet starting at row 905
# (dates 2019-11-18 .. 2019-11-28 that were previously missing), shifting
# the existing rows 905..977 down to 914..986.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 9 blank rows at row 905, shifting existing data down.
$insertRange = $ws.Range("A905:I913")
$insertRange.Insert(-4121)  # -4121 = xlShiftDown

# New rows data: timestamp, date, id, name, open, high, low, close, vol
$newRows = @(
    @(1574035200, "2019-11-18", "5273", "CHINHIN", 0.895, 0.915, 0.875, 0.875, 2206300),
    @(1574121600, "2019-11-19", "5273", "CHINHIN", 0.875, 0.905, 0.87,  0.885, 859400),
    @(1574208000, "2019-11-20", "5273", "CHINHIN", 0.885, 0.89,  0.88,  0.885, 417000),
    @(1574294400, "2019-11-21", "5273", "CHINHIN", 0.89,  0.905, 0.87,  0.895, 1444900),
    @(1574380800, "2019-11-22", "5273", "CHINHIN", 0.9,   0.92,  0.895, 0.905, 3186000),
    @(1574640000, "2019-11-25", "5273", "CHINHIN", 0.91,  0.915, 0.88,  0.9,   3123300),
    @(1574726400, "2019-11-26", "5273", "CHINHIN", 0.91,  0.92,  0.82,  0.825, 6263900),
    @(1574812800, "2019-11-27", "5273", "CHINHIN", 0.82,  0.835, 0.795, 0.8,   1485000),
    @(1574899200, "2019-11-28", "5273", "CHINHIN", 0.8,   0.8149999999999999, 0.795, 0.8, 1373500)
)

$startRow = 905
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $row[0]
    # Force the date and id columns to be stored as text (matching the
    # rest of the column), not auto-converted to a date serial / number.
    $ws.Cells.Item($r, 2).Value = "'" + $row[1]
    $ws.Cells.Item($r, 3).Value = "'" + $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = $row[8]
}
